# Auto-generated edit script: updates the cryptos price/volume table
# to match the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.734.38"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "2.528.71"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.61%  "

$ws.Range("E7").Value = "  -1.12%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").Value = "2.915.14"
$ws.Range("E14").Value = "  -1.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.59%  "

$ws.Range("D16").Value = "2.533.72"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("E17").Value = "  -3.87%  "

$ws.Range("D18").Value = "42.705.96"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("E24").Value = "  -3.06%  "

$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.33%  "

$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("E33").Value = "  +9.09%  "

$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0784"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.88%  "

$ws.Range("E38").Value = "  -7.17%  "

$ws.Range("E39").Value = "  -1.24%  "

$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.08%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").Value = "2.003.11"
$ws.Range("E46").Value = "  +0.24%  "

$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("D48").Value = "2.770.30"
$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.190"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.08%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.53%  "
